$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.122.93"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.781.28"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'" + "225.70"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").Value = "'" + "0.546"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'" + "32.10"
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("D9").Value = "'" + "0.293"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").Value = "'" + "0.0949"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "2.037.97"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "1.786.55"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "'" + "10.91"
$ws.Range("E14").Value = "  -5.08%  "
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "34.103.74"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "'" + "67.56"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "'" + "245.40"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "'" + "10.88"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "'" + "2.04"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").Value = "'" + "162.14"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").Value = "'" + "16.26"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").Value = "'" + "0.0517"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").Value = "'" + "3.72"
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("D35").Value = "1.446.60"
$ws.Range("E35").Value = "  +2.81%  "
$ws.Range("E36").Value = "  +5.68%  "
$ws.Range("D37").Value = "'" + "0.652"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").Value = "'" + "0.0191"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "'" + "81.25"
$ws.Range("E40").Value = "  +1.40%  "
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").Value = "'" + "0.913"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").Value = "'" + "13.62"
$ws.Range("E44").Value = "  +1.96%  "
$ws.Range("D45").Value = "'" + "0.0521"
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").Value = "1.938.19"
$ws.Range("D49").Value = "'" + "0.0" + [char]0x2086 + "0132"
$ws.Range("E49").Value = "  -6.81%  "
$ws.Range("D50").Value = "'" + "104.69"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("E51").Value = "  +0.31%  "

Write-Host "Applied 78 cell updates"